# update data prep for variable efficiency
#
# 1) "Definition" sheet: rows 7-18 (col A, node names) get reordered.
# 2) "Nodes" sheet: rows 2-13 (full rows, cols A-G) get reordered the same way.
# 3) "Object__to_from_node" sheet: new row 19 is appended
#    (unit__to_node / unit / Electrolyzer / Hydrogen_Kasso / ordered_unit_flow_op / TRUE)
#
# Row moves are done with Range.Cut(destination) via a scratch area far below
# the used range (and back) rather than re-typing values through .Value, so
# that cell data types (e.g. the literal text "true" vs. a native boolean)
# survive the reshuffle exactly as they were.

$wb = $excel.ActiveWorkbook

# Node names, in their original (before-edit) row order for rows 7..18 / 2..13.
$oldOrder = @(
    "Power_Kasso",
    "E-Methanol_Kasso",
    "E-Methanol_storage_Kasso",
    "Vaporized_Carbon_Dioxide",
    "Waste_Heat",
    "Carbon_Dioxide",
    "Hydrogen_Kasso",
    "Raw_Methanol",
    "District_Heating",
    "Water",
    "Hydrogen_storage_Kasso",
    "Power_Wholesale"
)

# The new (target) row order.
$newOrder = @(
    "Water",
    "Vaporized_Carbon_Dioxide",
    "Carbon_Dioxide",
    "Waste_Heat",
    "E-Methanol_Kasso",
    "Power_Wholesale",
    "E-Methanol_storage_Kasso",
    "District_Heating",
    "Hydrogen_Kasso",
    "Hydrogen_storage_Kasso",
    "Power_Kasso",
    "Raw_Methanol"
)

$rowCount = $newOrder.Length
$scratchBase = 500

# ---------------------------------------------------------------------------
# 1) Definition sheet - reorder column A for rows 7..18
# ---------------------------------------------------------------------------
$wsDef = $wb.Worksheets.Item("Definition")
$defStartRow = 7

for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $defStartRow + $i
    $dstRow = $scratchBase + $i
    $wsDef.Range("A$srcRow").Cut($wsDef.Range("A$dstRow"))
}
for ($i = 0; $i -lt $rowCount; $i++) {
    $name = $newOrder[$i]
    $oldIdx = [Array]::IndexOf($oldOrder, $name)
    $srcRow = $scratchBase + $oldIdx
    $dstRow = $defStartRow + $i
    $wsDef.Range("A$srcRow").Cut($wsDef.Range("A$dstRow"))
}

# ---------------------------------------------------------------------------
# 2) Nodes sheet - reorder full rows (A..G) for rows 2..13
# ---------------------------------------------------------------------------
$wsNodes = $wb.Worksheets.Item("Nodes")
$nodesStartRow = 2

for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $nodesStartRow + $i
    $dstRow = $scratchBase + $i
    $wsNodes.Range("A$srcRow" + ":G$srcRow").Cut($wsNodes.Range("A$dstRow" + ":G$dstRow"))
}
for ($i = 0; $i -lt $rowCount; $i++) {
    $name = $newOrder[$i]
    $oldIdx = [Array]::IndexOf($oldOrder, $name)
    $srcRow = $scratchBase + $oldIdx
    $dstRow = $nodesStartRow + $i
    $wsNodes.Range("A$srcRow" + ":G$srcRow").Cut($wsNodes.Range("A$dstRow" + ":G$dstRow"))
}

# ---------------------------------------------------------------------------
# 3) Object__to_from_node sheet - append new row 19
# ---------------------------------------------------------------------------
$wsRel = $wb.Worksheets.Item("Object__to_from_node")

$newRow = 19
$wsRel.Cells.Item($newRow, 1).Value = "unit__to_node"
$wsRel.Cells.Item($newRow, 2).Value = "unit"
$wsRel.Cells.Item($newRow, 3).Value = "Electrolyzer"
$wsRel.Cells.Item($newRow, 4).Value = "Hydrogen_Kasso"
$wsRel.Cells.Item($newRow, 5).Value = "ordered_unit_flow_op"
$wsRel.Cells.Item($newRow, 6).Value = $true

Write-Host "edit complete"
